# error solve ifrs list
# Rewrites the per-year financial figures on the (only) worksheet so the
# raw OOXML numbers match the corrected IFRS data set, and clears the
# now-unavailable trailing rows (7-9) down to just their label columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : 2014/12 ---------------------------------------------------
$ws.Range("D2").Value  = 3035
$ws.Range("E2").Value  = 159
$ws.Range("F2").Value  = 159
$ws.Range("G2").Value  = 271
$ws.Range("H2").Value  = 203
$ws.Range("I2").Value  = 167
$ws.Range("J2").Value  = 36
$ws.Range("K2").Value  = 3777
$ws.Range("L2").Value  = 1104
$ws.Range("M2").Value  = 2673
$ws.Range("N2").Value  = 2457
$ws.Range("O2").Value  = 216
$ws.Range("P2").Value  = 130
$ws.Range("Q2").Value  = 405
$ws.Range("R2").Value  = -307
$ws.Range("S2").Value  = -116
$ws.Range("T2").Value  = 82
$ws.Range("U2").Value  = 323
$ws.Range("V2").Value  = 102
$ws.Range("W2").Value  = 5.23
$ws.Range("X2").Value  = 6.7
$ws.Range("Y2").Value  = 7.01
$ws.Range("Z2").Value  = 5.56
$ws.Range("AA2").Value = 41.32
$ws.Range("AB2").Value = 1717.54
$ws.Range("AC2").Value = 643
$ws.Range("AD2").Value = 8.109999999999999
$ws.Range("AE2").Value = 9580
$ws.Range("AF2").Value = 0.54
$ws.Range("AG2").Value = 130
$ws.Range("AH2").Value = 2.49
$ws.Range("AI2").Value = 19.97
$ws.Range("AJ2").Value = 25947500

# --- Row 3 : 2015/12 ---------------------------------------------------
$ws.Range("D3").Value  = 2679
$ws.Range("E3").Value  = 50
$ws.Range("F3").Value  = 50
$ws.Range("G3").Value  = 224
$ws.Range("H3").Value  = 166
$ws.Range("I3").Value  = 120
$ws.Range("J3").Value  = 46
$ws.Range("K3").Value  = 3738
$ws.Range("L3").Value  = 912
$ws.Range("M3").Value  = 2826
$ws.Range("N3").Value  = 2569
$ws.Range("O3").Value  = 257
$ws.Range("P3").Value  = 130
$ws.Range("Q3").Value  = 450
$ws.Range("R3").Value  = -268
$ws.Range("S3").Value  = -108
$ws.Range("T3").Value  = 255
$ws.Range("U3").Value  = 195
$ws.Range("V3").Value  = 33
$ws.Range("W3").Value  = 1.85
$ws.Range("X3").Value  = 6.19
$ws.Range("Y3").Value  = 4.77
$ws.Range("Z3").Value  = 4.41
$ws.Range("AA3").Value = 32.28
$ws.Range("AB3").Value = 1790.61
$ws.Range("AC3").Value = 462
$ws.Range("AD3").Value = 8.99
$ws.Range("AE3").Value = 10017
$ws.Range("AF3").Value = 0.41
$ws.Range("AG3").Value = 120
$ws.Range("AH3").Value = 2.89
$ws.Range("AI3").Value = 25.69
$ws.Range("AJ3").Value = 25947500

# --- Row 4 : 2016/12 ---------------------------------------------------
$ws.Range("D4").Value  = 2578
$ws.Range("E4").Value  = 65
$ws.Range("F4").Value  = 65
$ws.Range("G4").Value  = 149
$ws.Range("H4").Value  = 108
$ws.Range("I4").Value  = 88
$ws.Range("J4").Value  = 20
$ws.Range("K4").Value  = 3773
$ws.Range("L4").Value  = 882
$ws.Range("M4").Value  = 2891
$ws.Range("N4").Value  = 2623
$ws.Range("O4").Value  = 268
$ws.Range("P4").Value  = 130
$ws.Range("Q4").Value  = 315
$ws.Range("R4").Value  = -300
$ws.Range("S4").Value  = -72
$ws.Range("T4").Value  = 102
$ws.Range("U4").Value  = 213
$ws.Range("V4").Value  = $null
$ws.Range("W4").Value  = 2.53
$ws.Range("X4").Value  = 4.19
$ws.Range("Y4").Value  = 3.37
$ws.Range("Z4").Value  = 2.88
$ws.Range("AA4").Value = 30.5
$ws.Range("AB4").Value = 1841.06
$ws.Range("AC4").Value = 338
$ws.Range("AD4").Value = 13.14
$ws.Range("AE4").Value = 10226
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 120
$ws.Range("AH4").Value = 2.71
$ws.Range("AI4").Value = 35.14
$ws.Range("AJ4").Value = 25947500

# --- Row 5 : 2017/12 ---------------------------------------------------
$ws.Range("D5").Value  = 2517
$ws.Range("E5").Value  = 32
$ws.Range("F5").Value  = 32
$ws.Range("G5").Value  = 102
$ws.Range("H5").Value  = 67
$ws.Range("I5").Value  = 62
$ws.Range("J5").Value  = 5
$ws.Range("K5").Value  = 3827
$ws.Range("L5").Value  = 910
$ws.Range("M5").Value  = 2918
$ws.Range("N5").Value  = 2650
$ws.Range("O5").Value  = 267
$ws.Range("P5").Value  = 130
$ws.Range("Q5").Value  = 209
$ws.Range("R5").Value  = -136
$ws.Range("S5").Value  = -38
$ws.Range("T5").Value  = 51
$ws.Range("U5").Value  = 158
$ws.Range("V5").Value  = 0
$ws.Range("W5").Value  = 1.29
$ws.Range("X5").Value  = 2.67
$ws.Range("Y5").Value  = 2.37
$ws.Range("Z5").Value  = 1.77
$ws.Range("AA5").Value = 31.18
$ws.Range("AB5").Value = 1871.64
$ws.Range("AC5").Value = 240
$ws.Range("AD5").Value = 14.81
$ws.Range("AE5").Value = 10334
$ws.Range("AF5").Value = 0.34
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 3.37
$ws.Range("AI5").Value = 49.36
$ws.Range("AJ5").Value = 25947500

# --- Row 6 : 2018/12 (note: J6/O6 already absent in the source row) ---
$ws.Range("D6").Value  = 2466
$ws.Range("E6").Value  = 18
$ws.Range("F6").Value  = 18
$ws.Range("G6").Value  = 89
$ws.Range("H6").Value  = 51
$ws.Range("I6").Value  = 37
$ws.Range("K6").Value  = 3756
$ws.Range("L6").Value  = 868
$ws.Range("M6").Value  = 2888
$ws.Range("N6").Value  = 2614
$ws.Range("P6").Value  = 130
$ws.Range("Q6").Value  = 104
$ws.Range("R6").Value  = -25
$ws.Range("S6").Value  = -36
$ws.Range("T6").Value  = 92
$ws.Range("U6").Value  = 12
$ws.Range("V6").Value  = 0
$ws.Range("W6").Value  = 0.73
$ws.Range("X6").Value  = 2.07
$ws.Range("Y6").Value  = 1.4
$ws.Range("Z6").Value  = 1.34
$ws.Range("AA6").Value = 30.07
$ws.Range("AB6").Value = 1867.95
$ws.Range("AC6").Value = 142
$ws.Range("AD6").Value = 20.53
$ws.Range("AE6").Value = 10193
$ws.Range("AF6").Value = 0.29
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 3.44
$ws.Range("AI6").Value = 69.73
$ws.Range("AJ6").Value = 25947500

# --- Rows 7-9 : estimate years (2019/12(E)-2021/12(E)) no longer carry
#     any of the projected figures -- only the label columns (A-C) stay.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
